$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Build the full refreshed table (header + 39 holdings rows) as a 2D array,
# then write it in one shot to A1:H40.
$arr = New-Object 'object[,]' 40,8

$arr[0,0] = "ISIN"
$arr[0,1] = "Stock Name"
$arr[0,2] = "Mutual Fund"
$arr[0,3] = "Jan_2026"
$arr[0,4] = "Dec_2025"
$arr[0,5] = "Nov_2025"
$arr[0,6] = "MoM"
$arr[0,7] = "QoQ"

$arr[1,0] = "INE040A01034"
$arr[1,1] = "HDFC Bank Limited"
$arr[1,2] = "quant Multi Asset Allocation Fund"
$arr[1,3] = 9.305626
$arr[1,4] = 8.382281000000001
$arr[1,5] = 0
$arr[1,6] = 0.9233449999999994
$arr[1,7] = 9.305626

$arr[2,0] = "INE090A01021"
$arr[2,1] = "ICICI Bank Limited"
$arr[2,2] = "quant Multi Asset Allocation Fund"
$arr[2,3] = 9.26423
$arr[2,4] = 8.473186999999999
$arr[2,5] = 2.94592
$arr[2,6] = 0.7910430000000002
$arr[2,7] = 6.318309999999999

$arr[3,0] = "INE237A01036"
$arr[3,1] = "Kotak Mahindra Bank Limited"
$arr[3,2] = "quant Multi Asset Allocation Fund"
$arr[3,3] = 7.354324
$arr[3,4] = 0
$arr[3,5] = 0
$arr[3,6] = 7.354324
$arr[3,7] = 7.354324

$arr[4,0] = "INE795G01014"
$arr[4,1] = "HDFC Life Insurance Co Ltd"
$arr[4,2] = "quant Multi Asset Allocation Fund"
$arr[4,3] = 5.029332
$arr[4,4] = 5.541729
$arr[4,5] = 3.065313
$arr[4,6] = -0.512397
$arr[4,7] = 1.964019

$arr[5,0] = "INE296A01032"
$arr[5,1] = "Bajaj Finance Limited"
$arr[5,2] = "quant Multi Asset Allocation Fund"
$arr[5,3] = 4.362969
$arr[5,4] = 4.973663
$arr[5,5] = 0
$arr[5,6] = -0.6106940000000005
$arr[5,7] = 4.362969

$arr[6,0] = "INE075A01022"
$arr[6,1] = "Wipro Ltd"
$arr[6,2] = "quant Multi Asset Allocation Fund"
$arr[6,3] = 4.315831
$arr[6,4] = 0
$arr[6,5] = 0
$arr[6,6] = 4.315831
$arr[6,7] = 4.315831

$arr[7,0] = "INE0BS701011"
$arr[7,1] = "Premier Energies Limited"
$arr[7,2] = "quant Multi Asset Allocation Fund"
$arr[7,3] = 3.362745
$arr[7,4] = 4.225607
$arr[7,5] = 5.191142
$arr[7,6] = -0.8628620000000002
$arr[7,7] = -1.828397

$arr[8,0] = "INE406A01037"
$arr[8,1] = "Aurobindo Pharma Limited"
$arr[8,2] = "quant Multi Asset Allocation Fund"
$arr[8,3] = 2.966447
$arr[8,4] = 0
$arr[8,5] = 0
$arr[8,6] = 2.966447
$arr[8,7] = 2.966447

$arr[9,0] = "INE127D01025"
$arr[9,1] = "HDFC Asset Management Company Ltd"
$arr[9,2] = "quant Multi Asset Allocation Fund"
$arr[9,3] = 2.387782
$arr[9,4] = 2.724474
$arr[9,5] = 0
$arr[9,6] = -0.3366919999999998
$arr[9,7] = 2.387782

$arr[10,0] = "INE364U01010"
$arr[10,1] = "Adani Green Energy Limited"
$arr[10,2] = "quant Multi Asset Allocation Fund"
$arr[10,3] = 2.231143
$arr[10,4] = 2.854445
$arr[10,5] = 0
$arr[10,6] = -0.6233020000000002
$arr[10,7] = 2.231143

$arr[11,0] = "INE018A01030"
$arr[11,1] = "Larsen & Toubro Limited"
$arr[11,2] = "quant Multi Asset Allocation Fund"
$arr[11,3] = 2.128037
$arr[11,4] = 1.476168
$arr[11,5] = 0
$arr[11,6] = 0.651869
$arr[11,7] = 2.128037

$arr[12,0] = "INE261F16AE9"
$arr[12,1] = "NABARD CD 19-Jan-2027"
$arr[12,2] = "quant Multi Asset Allocation Fund"
$arr[12,3] = 1.962753
$arr[12,4] = 0
$arr[12,5] = 0
$arr[12,6] = 1.962753
$arr[12,7] = 1.962753

$arr[13,0] = "INE154A01025"
$arr[13,1] = "ITC Limited"
$arr[13,2] = "quant Multi Asset Allocation Fund"
$arr[13,3] = 1.932332
$arr[13,4] = 2.59661
$arr[13,5] = 2.76184
$arr[13,6] = -0.6642780000000001
$arr[13,7] = -0.8295079999999999

$arr[14,0] = "INE081A01020"
$arr[14,1] = "Tata Steel Limited"
$arr[14,2] = "quant Multi Asset Allocation Fund"
$arr[14,3] = 1.552077
$arr[14,4] = 0.7772790000000001
$arr[14,5] = 0
$arr[14,6] = 0.7747979999999999
$arr[14,7] = 1.552077

$arr[15,0] = "INE676A01027"
$arr[15,1] = "Black Box Limited"
$arr[15,2] = "quant Multi Asset Allocation Fund"
$arr[15,3] = 1.530466
$arr[15,4] = 1.757758
$arr[15,5] = 1.806139
$arr[15,6] = -0.2272919999999998
$arr[15,7] = -0.2756729999999998

$arr[16,0] = "INE514E16CN1"
$arr[16,1] = "EXIM Bank CD 11-Nov-2026"
$arr[16,2] = "quant Multi Asset Allocation Fund"
$arr[16,3] = 1.491614
$arr[16,4] = 1.600926
$arr[16,5] = 0
$arr[16,6] = -0.1093120000000001
$arr[16,7] = 1.491614

$arr[17,0] = "INE397D01024"
$arr[17,1] = "Bharti Airtel Limited"
$arr[17,2] = "quant Multi Asset Allocation Fund"
$arr[17,3] = 1.450806
$arr[17,4] = 0.780396
$arr[17,5] = 0
$arr[17,6] = 0.6704100000000001
$arr[17,7] = 1.450806

$arr[18,0] = "INE155A01022"
$arr[18,1] = "Tata Motors Passenger Vehicles Limited"
$arr[18,2] = "quant Multi Asset Allocation Fund"
$arr[18,3] = 1.425685
$arr[18,4] = 0.732985
$arr[18,5] = 0
$arr[18,6] = 0.6927000000000001
$arr[18,7] = 1.425685

$arr[19,0] = "INE781S01027"
$arr[19,1] = "Ventive Hospitality Limited"
$arr[19,2] = "quant Multi Asset Allocation Fund"
$arr[19,3] = 1.242275
$arr[19,4] = 1.364569
$arr[19,5] = 1.381006
$arr[19,6] = -0.1222939999999999
$arr[19,7] = -0.1387309999999999

$arr[20,0] = "INE271C01023"
$arr[20,1] = "DLF Limited"
$arr[20,2] = "quant Multi Asset Allocation Fund"
$arr[20,3] = 1.070117
$arr[20,4] = 1.242889
$arr[20,5] = 1.387293
$arr[20,6] = -0.1727719999999999
$arr[20,7] = -0.3171760000000001

$arr[21,0] = "INE205A01025"
$arr[21,1] = "Vedanta Limited"
$arr[21,2] = "quant Multi Asset Allocation Fund"
$arr[21,3] = 0.809566
$arr[21,4] = 0.771182
$arr[21,5] = 0
$arr[21,6] = 0.03838399999999997
$arr[21,7] = 0.809566

$arr[22,0] = "INE918I01026"
$arr[22,1] = "Bajaj Finserv Ltd."
$arr[22,2] = "quant Multi Asset Allocation Fund"
$arr[22,3] = 0.755288
$arr[22,4] = 0.793538
$arr[22,5] = 0
$arr[22,6] = -0.03825000000000001
$arr[22,7] = 0.755288

$arr[23,0] = "INE200M01039"
$arr[23,1] = "Varun Beverages Limited"
$arr[23,2] = "quant Multi Asset Allocation Fund"
$arr[23,3] = 0.700078
$arr[23,4] = 0.781693
$arr[23,5] = 0
$arr[23,6] = -0.08161499999999999
$arr[23,7] = 0.700078

$arr[24,0] = "INE376G01013"
$arr[24,1] = "Biocon Ltd"
$arr[24,2] = "quant Multi Asset Allocation Fund"
$arr[24,3] = 0.64877
$arr[24,4] = 0.748386
$arr[24,5] = 0
$arr[24,6] = -0.09961600000000004
$arr[24,7] = 0.64877

$arr[25,0] = "INE1JAR25012"
$arr[25,1] = "Knowledge Realty Trust"
$arr[25,2] = "quant Multi Asset Allocation Fund"
$arr[25,3] = 0.644485
$arr[25,4] = 0.68784
$arr[25,5] = 0.711893
$arr[25,6] = -0.04335500000000003
$arr[25,7] = -0.06740800000000002

$arr[26,0] = "INE484J01027"
$arr[26,1] = "Godrej Properties Limited"
$arr[26,2] = "quant Multi Asset Allocation Fund"
$arr[26,3] = 0.602641
$arr[26,4] = 0.822896
$arr[26,5] = 0
$arr[26,6] = -0.220255
$arr[26,7] = 0.602641

$arr[27,0] = "INE549H01021"
$arr[27,1] = "Anand Rathi Share & Stock Brokers Ltd"
$arr[27,2] = "quant Multi Asset Allocation Fund"
$arr[27,3] = 0.58982
$arr[27,4] = 0.689252
$arr[27,5] = 0.804961
$arr[27,6] = -0.09943199999999996
$arr[27,7] = -0.215141

$arr[28,0] = "INE752E01010"
$arr[28,1] = "Power Grid Corporation of India Limited"
$arr[28,2] = "quant Multi Asset Allocation Fund"
$arr[28,3] = 0.541238
$arr[28,4] = 0.5997479999999999
$arr[28,5] = 0
$arr[28,6] = -0.05850999999999995
$arr[28,7] = 0.541238

$arr[29,0] = "INE414G14UT3"
$arr[29,1] = "Muthoot Finance Ltd CP 08-Sep-2026"
$arr[29,2] = "quant Multi Asset Allocation Fund"
$arr[29,3] = 0.300465
$arr[29,4] = 0.321699
$arr[29,5] = 0.339716
$arr[29,6] = -0.02123400000000003
$arr[29,7] = -0.03925100000000004

$arr[30,0] = "INE556F16BN1"
$arr[30,1] = "SIDBI CD 27-Oct-2026"
$arr[30,2] = "quant Multi Asset Allocation Fund"
$arr[30,3] = 0.199234
$arr[30,4] = 0.213812
$arr[30,5] = 0.225869
$arr[30,6] = -0.01457800000000001
$arr[30,7] = -0.02663499999999999

$arr[31,0] = "INE062A01020"
$arr[31,1] = "State Bank of India"
$arr[31,2] = "quant Multi Asset Allocation Fund"
$arr[31,3] = 0
$arr[31,4] = 0
$arr[31,5] = 7.226046
$arr[31,6] = 0
$arr[31,7] = -7.226046

$arr[32,0] = "INE0NHL23019"
$arr[32,1] = "Indus Infra Trust"
$arr[32,2] = "quant Multi Asset Allocation Fund"
$arr[32,3] = 0
$arr[32,4] = 0
$arr[32,5] = 2.480624
$arr[32,6] = 0
$arr[32,7] = -2.480624

$arr[33,0] = "INE758E01017"
$arr[33,1] = "Jio Financial Services Limited"
$arr[33,2] = "quant Multi Asset Allocation Fund"
$arr[33,3] = 0
$arr[33,4] = 3.345709
$arr[33,5] = 3.682917
$arr[33,6] = -3.345709
$arr[33,7] = -3.682917

$arr[34,0] = "INE245A01021"
$arr[34,1] = "Tata Power Company Limited"
$arr[34,2] = "quant Multi Asset Allocation Fund"
$arr[34,3] = 0
$arr[34,4] = 0.58974
$arr[34,5] = 2.374302
$arr[34,6] = -0.58974
$arr[34,7] = -2.374302

$arr[35,0] = "INE0Z8Z23013"
$arr[35,1] = "Capital Infra Trust InvIT"
$arr[35,2] = "quant Multi Asset Allocation Fund"
$arr[35,3] = 0
$arr[35,4] = 0
$arr[35,5] = 0.921096
$arr[35,6] = 0
$arr[35,7] = -0.921096

$arr[36,0] = "INE556F16AY0"
$arr[36,1] = "SIDBI CD 13-Jan-2026"
$arr[36,2] = "quant Multi Asset Allocation Fund"
$arr[36,3] = 0
$arr[36,4] = 2.250632
$arr[36,5] = 2.37337
$arr[36,6] = -2.250632
$arr[36,7] = -2.37337

$arr[37,0] = "INE467B01029"
$arr[37,1] = "Tata Consultancy Services Limited"
$arr[37,2] = "quant Multi Asset Allocation Fund"
$arr[37,3] = 0
$arr[37,4] = 0
$arr[37,5] = 2.863599
$arr[37,6] = 0
$arr[37,7] = -2.863599

$arr[38,0] = "INE237A01028"
$arr[38,1] = "Kotak Mahindra Bank Limited"
$arr[38,2] = "quant Multi Asset Allocation Fund"
$arr[38,3] = 0
$arr[38,4] = 8.523743
$arr[38,5] = 0
$arr[38,6] = -8.523743
$arr[38,7] = 0

$arr[39,0] = "INE002A01018"
$arr[39,1] = "Reliance Industries Limited"
$arr[39,2] = "quant Multi Asset Allocation Fund"
$arr[39,3] = 0
$arr[39,4] = 0.148326
$arr[39,5] = 0
$arr[39,6] = -0.148326
$arr[39,7] = 0

$ws.Range("A1:H40").Value = $arr